$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.022.26"
$ws.Range("E2").Value = "  +0.09%  "
$ws.Range("D3").Value = "3.877.31"
$ws.Range("E3").Value = "  -0.94%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "482.60"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.20%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.05"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.08%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.620"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.736"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.59%  "
$ws.Range("E10").Value = "  +7.08%  "
$ws.Range("E11").Value = "  -0.37%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "42.71"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.60%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.53"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.36%  "
$ws.Range("D14").Value = "4.495.07"
$ws.Range("E14").Value = "  -0.77%  "
$ws.Range("D15").Value = "3.871.72"
$ws.Range("E15").Value = "  -3.48%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.28"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.59%  "
$ws.Range("E17").Value = "  -0.69%  "
$ws.Range("E18").Value = "  +1.39%  "
$ws.Range("E19").Value = "  -0.38%  "
$ws.Range("D20").Value = "68.056.40"
$ws.Range("E20").Value = "  +0.10%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "428.64"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.83%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.54"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.50%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.73"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.95%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "90.07"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.54%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.96"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +10.07%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.68"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.37%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.99"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.29%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "37.37"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.27%  "
$ws.Range("E29").Value = "  -3.97%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "711.60"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.22%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "13.49"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.57%  "
$ws.Range("E32").Value = "  +0.13%  "
$ws.Range("E33").Value = "  +2.80%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.06"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +9.94%  "
$ws.Range("D35").Value = "0.0₃0871"
$ws.Range("E35").Value = "  -2.52%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "40.88"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.12%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "60.80"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.98%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0500"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.57%  "
$ws.Range("B39").Value = "Dai"
$ws.Range("C39").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.998"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.03%  "
$ws.Range("E40").Value = "  -3.97%  "
$ws.Range("B41").Value = "TheGraph"
$ws.Range("C41").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.395"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +13.10%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.97"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.29%  "
$ws.Range("E43").Value = "  +3.52%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.98"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.17%  "
$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.142"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.17%  "
$ws.Range("B46").Value = "ApeXProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.36"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.99%  "
$ws.Range("E47").Value = "  +0.18%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.38"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.04%  "
$ws.Range("E49").Value = "  -3.05%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "144.94"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.59%  "
$ws.Range("E51").Value = "  -1.45%  "
